$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.687.18'
$ws.Range('D3').Value = '1.635.15'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''212.87'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  +2.13%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').Value = '''19.01'
$ws.Range('E10').Value = '  +2.99%  '
$ws.Range('E11').Value = '  +2.72%  '
$ws.Range('D12').Value = '1.863.27'
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('D13').Value = '1.640.65'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').Value = '26.689.15'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('D17').Value = '''63.03'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '''208.56'
$ws.Range('E20').Value = '  +3.82%  '
$ws.Range('D21').Value = '''4.31'
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('D22').Value = '''9.40'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('E23').Value = '  +2.87%  '
$ws.Range('D24').Value = '''1.91'
$ws.Range('E24').Value = '  +2.34%  '
$ws.Range('D25').Value = '''146.49'
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = '''6.74'
$ws.Range('E28').Value = '  +2.61%  '
$ws.Range('D29').Value = '''15.39'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('E30').Value = '  +5.52%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('E34').Value = '  +1.51%  '
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('D36').Value = '1.168.81'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').Value = '''0.807'
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '''0.504'
$ws.Range('E40').Value = '  +1.47%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '''2.32'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').Value = '''0.795'
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('D43').Value = '''5.38'
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D44').Value = '1.773.87'
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('D45').Value = '''92.43'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').Value = '''54.70'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').Value = '''0.409'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '''7.54'
$ws.Range('E51').Value = '  +4.56%  '

# Restore default (unstyled) style for cells forced to text via apostrophe prefix
$ws.Range('D5').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
